# Adds a new "Date added in KEGG" column (E) to the KEGG immune-system
# table, fills in dates (or "NA" where unknown) for every pathway, and
# highlights the two rows that are discussed in the figures (fig 1 / fig 2
# per the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell -------------------------------------------------
$ws.Range("E1").Value = "Date added in KEGG"

# --- New column E data -------------------------------------------------
# "NA" rows (date unknown / not applicable)
$naRows = 2,3,7,12,14,18,21
foreach ($r in $naRows) {
    $ws.Cells.Item($r, 5).Value = "NA"
}

# Rows with a known "date added in KEGG", expressed as Excel date serials
# (days since 1899-12-30) so there is no locale-dependent parsing:
#   row -> serial (date)
#   4   -> 41757 (2014-04-28)
#   5   -> 44214 (2021-01-18)
#   6   -> 42613 (2016-08-31)
#   8   -> 40113 (2009-10-27)
#   9   -> 40072 (2009-09-16)
#   10  -> 40146 (2009-11-29)
#   11  -> 43185 (2018-03-26)
#   13  -> 38756 (2006-02-08)
#   15  -> 42688 (2016-11-14)
#   16  -> 42718 (2016-12-14)
#   17  -> 42746 (2017-01-11)
#   19  -> 38775 (2006-02-27)
#   20  -> 38953 (2006-08-24)
#   22  -> 40172 (2009-12-25)
#   23  -> 38885 (2006-06-17)
$dateSerials = @{
    4  = 41757
    5  = 44214
    6  = 42613
    8  = 40113
    9  = 40072
    10 = 40146
    11 = 43185
    13 = 38756
    15 = 42688
    16 = 42718
    17 = 42746
    19 = 38775
    20 = 38953
    22 = 40172
    23 = 38885
}

foreach ($r in $dateSerials.Keys) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $dateSerials[$r]
    $cell.NumberFormat = "YYYY\-MM\-DD"
}

# --- Highlight the two rows referenced in the figures -----------------
$highlightColor = 10420223   # RGB(255,255,158) light-yellow highlight
foreach ($r in 5,11) {
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 5)).Interior.Color = $highlightColor
}

# --- Column E width ------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 9.14

# --- Sheet view: zoom + selection ------------------------------------
$excel.ActiveWindow.Zoom = 110
[void]$ws.Range("H12").Select()
